# "Add files via upload" — re-save of Book1.xlsx where a batch of rows in
# the "Type" (genre) column, C, were re-typed from lower-case shared
# strings ("pop", "pop rock", "indie pop", "country pop", "dance pop",
# "soul pop", "acoustic pop", "folk") to capitalized ones ("Pop", "Pop
# rock", ...). Because the new text doesn't match any existing shared
# string, writing it via .Value appends a fresh entry to the shared
# string table (uniqueCount 183 -> 191) instead of reusing the old index,
# exactly like the diff shows (e.g. C32 6 -> 183).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C32").Value = "Pop"
$ws.Range("C33").Value = "Pop"
$ws.Range("C34").Value = "Pop rock"
$ws.Range("C38").Value = "Pop"
$ws.Range("C41").Value = "Pop"
$ws.Range("C43").Value = "Pop"
$ws.Range("C44").Value = "Indie pop"
$ws.Range("C45").Value = "Country pop"
$ws.Range("C46").Value = "Pop rock"
$ws.Range("C48").Value = "Dance pop"
$ws.Range("C49").Value = "Pop rock"
$ws.Range("C50").Value = "Dance pop"
$ws.Range("C52").Value = "Pop rock"
$ws.Range("C53").Value = "Pop rock"
$ws.Range("C55").Value = "Dance pop"
$ws.Range("C56").Value = "Pop rock"
$ws.Range("C59").Value = "Soul pop"
$ws.Range("C60").Value = "Pop rock"
$ws.Range("C61").Value = "Country pop"
$ws.Range("C63").Value = "Pop"
$ws.Range("C64").Value = "Pop rock"
$ws.Range("C65").Value = "Pop"
$ws.Range("C66").Value = "Acoustic pop"
$ws.Range("C67").Value = "Indie pop"
$ws.Range("C68").Value = "Pop"
$ws.Range("C69").Value = "Folk"
$ws.Range("C70").Value = "Dance pop"
$ws.Range("C71").Value = "Acoustic pop"
$ws.Range("C72").Value = "Dance pop"

# Match the author's final cursor position/selection and scroll before
# saving (sheet view had scrolled down so row 70 was the top visible row).
$ws.Range("C73").Select()
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
